# Refactor PPV Tools -- add "Program" and "Lot" columns to the `ppv`
# table on the PPV worksheet (update 20251124).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PPV")

# The PPV sheet holds a single Excel Table named "ppv" (A1:H2 before the
# edit). Grow it by two columns using the real ListObject/ListColumns
# object model, same as a user dragging the table's resize handle / typing
# new header text to the right of the table.
$lo = $ws.ListObjects.Item("ppv")

$colProgram = $lo.ListColumns.Add()
$colProgram.Range.Cells(1, 1).Value = "Program"

$colLot = $lo.ListColumns.Add()
$colLot.Range.Cells(1, 1).Value = "Lot"

# Match the look of the rest of the header row (bold white text on the
# table's banded header fill + the thin/blue-accent header border) by
# copying the adjacent header cell's formatting onto the two new headers.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Give the two new (currently empty) data cells the thin left/right column
# separators the rest of the table body uses.
foreach ($addr in @("I2", "J2")) {
    $cell = $ws.Range($addr)
    $cell.Borders.Item(7).LineStyle = 1    # xlEdgeLeft
    $cell.Borders.Item(7).Weight = 2       # xlThin
    $cell.Borders.Item(10).LineStyle = 1   # xlEdgeRight
    $cell.Borders.Item(10).Weight = 2      # xlThin
}

# Leave the selection where the edit ended, on the new last cell.
$ws.Range("J2").Select() | Out-Null
